$d = $word.ActiveDocument

# --- Change 1 ---------------------------------------------------------
# Remove the whole "Colour based on level of crawling " bullet paragraph
# (it sat directly above "Colour based on user activity on the articles...").
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Colour based on level of crawling ") {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}

# --- Change 2 ---------------------------------------------------------
# "Relevance" becomes "Filter", and the following bullet
# "Importance of article" is removed entirely.
$null = $d.Content.Find.Execute("Relevance", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "Filter", 2)

$target2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Importance of article") {
        $target2 = $p
        break
    }
}
if ($target2 -ne $null) {
    $target2.Range.Delete()
}

# --- Change 3 ---------------------------------------------------------
# The "_GoBack" bookmark moves from the end of the "Letter vs full name..."
# bullet to the middle of the word "Moving" (splitting it into "Mov" | "ing
# nodes") in the "Moving nodes" bullet.
$oldBm = $d.Bookmarks.Item("_GoBack")
$oldBm.Delete()

$movingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Moving nodes") {
        $movingPara = $p
        break
    }
}
$splitPos = $movingPara.Range.Start + 3
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
